# Trade #120 closed at 2026-02-17 16:04:05 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up figures to account for the
# newly-closed trade, and appends the trade's row to both the "All Trades"
# and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.84   # Current Capital
$summary.Range("B4").Value = -1.17    # Total P&L $
$summary.Range("B5").Value = -0.19    # Total P&L %
$summary.Range("B6").Value = 120      # Total Trades
$summary.Range("B7").Value = 44       # Winning Trades
$summary.Range("B9").Value = 36.67    # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.84     # Capital
$status.Range("D4").Value = 120       # Trades
$status.Range("E4").Value = -1.17     # P&L $
$status.Range("F4").Value = -1.16     # P&L %
$status.Range("G4").Value = 36.67     # Win Rate %

# ---------------------------------------------------------------------
# New trade row (#120) appended to "All Trades" and "MarketMaking" logs
# ---------------------------------------------------------------------
# Note: column B holds a date-like string ("2026-02-17"). Assigning it
# as a bare string lets Excel's input-parsing reinterpret it as a real
# date serial, which the source log never used (every other row stores
# it as literal text). A leading apostrophe forces text entry instead,
# matching the existing rows.
$targetRow = 121

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item($targetRow, 1).Value = 120
    $ws.Cells.Item($targetRow, 2).Value = "'2026-02-17"
    $ws.Cells.Item($targetRow, 3).Value = "16:03:58"
    $ws.Cells.Item($targetRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($targetRow, 5).Value = "UP"
    $ws.Cells.Item($targetRow, 6).Value = 0.91
    $ws.Cells.Item($targetRow, 7).Value = 0.92
    $ws.Cells.Item($targetRow, 8).Value = "CLOSED"
    $ws.Cells.Item($targetRow, 9).Value = 1.0989
    $ws.Cells.Item($targetRow, 10).Value = 0.01
    $ws.Cells.Item($targetRow, 11).Value = 98.84
    $ws.Cells.Item($targetRow, 12).Value = 0
    $ws.Cells.Item($targetRow, 13).Value = 0
    $ws.Cells.Item($targetRow, 14).Value = 0.6
    $ws.Cells.Item($targetRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($targetRow, 16).Value = "early_exit"
    $ws.Cells.Item($targetRow, 17).Value = 0.15
}
